# msz - more testcases, 2 new selectDropdown keywords, Videorecording
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 17 (104_MotorcycleInsurance_001_SmokeTest),
# pushing the existing rows 17 (Motorcycle) and 18 (Camper) down to 20 and 21.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# Fill column A (Record/Process) for the first two new rows
$ws.Range("A17").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields"
$ws.Range("A18").Value = "103_TruckInsurance_003_InsurantData_002_FieldHintsAndErrors"

# Fill column B (Variables) for the first two new rows
$ws.Range("B17").Value = "var103_TruckInsurance_003_InsurantData_001_MandatoryFields"
$ws.Range("B18").Value = "var103_TruckInsurance_003_InsurantData_002_FieldHintsAndErrors"

# Fill column A/B for the third new row
$ws.Range("A19").Value = "103_TruckInsurance_003_InsurantData_003_ListContents"
$ws.Range("B19").Value = "var103_TruckInsurance_003_InsurantData_003_ListContents"

# Fill column C (link to proTruckInsurance) for all three new rows
$ws.Range("C17").Value = "Open Truck Insurance"
$ws.Range("C18").Value = "Open Truck Insurance"
$ws.Range("C19").Value = "Open Truck Insurance"

# Fill column E (proTruckInsurance test marker) for all three new rows
$ws.Range("E17").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields"
$ws.Range("E18").Value = "103_TruckInsurance_003_InsurantData_002_FieldHintsAndErrors"
$ws.Range("E19").Value = "103_TruckInsurance_003_InsurantData_003_ListContents"

# Update selection to match the new authored state
$ws.Range("E34").Select()
